$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44364
$ws.Range("K2").Value = 13000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 13500
$ws.Range("P2").Value = 750
$ws.Range("D3").Value = 44364
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 11000
$ws.Range("M3").Value = 11000
$ws.Range("P3").Value = 611
$ws.Range("D4").Value = 44379
$ws.Range("D5").Value = 44379
$ws.Range("J5").Value = 100
$ws.Range("D6").Value = 44259
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 12500
$ws.Range("P6").Value = 694
$ws.Range("D7").Value = 44259
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 10000
$ws.Range("P7").Value = 556
$ws.Range("D8").Value = 44308
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 10000
$ws.Range("L8").Value = 11000
$ws.Range("M8").Value = 10500
$ws.Range("P8").Value = 583
$ws.Range("D9").Value = 44308
$ws.Range("D10").Value = 44342
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 11500
$ws.Range("P10").Value = 639
$ws.Range("D11").Value = 44342
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 50
$ws.Range("K11").Value = 9000
$ws.Range("L11").Value = 9000
$ws.Range("M11").Value = 9000
$ws.Range("P11").Value = 500
$ws.Range("D12").Value = 44253
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 12000
$ws.Range("P12").Value = 667
$ws.Range("D13").Value = 44253
$ws.Range("I13").Value = "Segunda"
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 10000
$ws.Range("P13").Value = 556
$ws.Range("D14").Value = 44384
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 16000
$ws.Range("M14").Value = 15500
$ws.Range("P14").Value = 861
$ws.Range("D15").Value = 44384
$ws.Range("J15").Value = 50
$ws.Range("K15").Value = 13000
$ws.Range("L15").Value = 13000
$ws.Range("M15").Value = 13000
$ws.Range("P15").Value = 722
$ws.Range("D18").Value = 44392
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 15000
$ws.Range("L18").Value = 16000
$ws.Range("M18").Value = 15500
$ws.Range("P18").Value = 861
$ws.Range("D19").Value = 44392
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 14000
$ws.Range("L19").Value = 14000
$ws.Range("M19").Value = 14000
$ws.Range("P19").Value = 778
$ws.Range("D20").Value = 44328
$ws.Range("J20").Value = 200
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = 9500
$ws.Range("P20").Value = 528
$ws.Range("D21").Value = 44328
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 8000
$ws.Range("L21").Value = 8000
$ws.Range("M21").Value = 8000
$ws.Range("P21").Value = 444
$ws.Range("D22").Value = 44265
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 13000
$ws.Range("L22").Value = 14000
$ws.Range("M22").Value = 13500
$ws.Range("P22").Value = 750
$ws.Range("D23").Value = 44316
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 11000
$ws.Range("M23").Value = 10500
$ws.Range("P23").Value = 583
$ws.Range("D24").Value = 44316
$ws.Range("I24").Value = "Segunda"
$ws.Range("K24").Value = 9000
$ws.Range("L24").Value = 9000
$ws.Range("M24").Value = 9000
$ws.Range("P24").Value = 500
$ws.Range("D25").Value = 44320
$ws.Range("I25").Value = "Primera"
$ws.Range("K25").Value = 9000
$ws.Range("M25").Value = 9500
$ws.Range("P25").Value = 528
$ws.Range("D26").Value = 44320
$ws.Range("I26").Value = "Segunda"
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 8000
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = 8000
$ws.Range("P26").Value = 444
$ws.Range("D27").Value = 44313
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 10000
$ws.Range("L27").Value = 11000
$ws.Range("M27").Value = 10500
$ws.Range("P27").Value = 583
$ws.Range("D28").Value = 44313
$ws.Range("I28").Value = "Segunda"
$ws.Range("K28").Value = 9000
$ws.Range("L28").Value = 9000
$ws.Range("M28").Value = 9000
$ws.Range("P28").Value = 500
$ws.Range("D29").Value = 44356
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 100
$ws.Range("L29").Value = 11000
$ws.Range("M29").Value = 10500
$ws.Range("P29").Value = 583
$ws.Range("D30").Value = 44356
$ws.Range("I30").Value = "Segunda"
$ws.Range("J30").Value = 50
$ws.Range("L30").Value = 9000
$ws.Range("M30").Value = 9000
$ws.Range("P30").Value = 500
$ws.Range("D31").Value = 44280
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 10000
$ws.Range("L31").Value = 11000
$ws.Range("M31").Value = 10500
$ws.Range("P31").Value = 583
$ws.Range("D32").Value = 44280
$ws.Range("I32").Value = "Segunda"
$ws.Range("K32").Value = 9000
$ws.Range("L32").Value = 9000
$ws.Range("M32").Value = 9000
$ws.Range("P32").Value = 500
$ws.Range("D33").Value = 44397
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 14000
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = 14500
$ws.Range("P33").Value = 806
$ws.Range("D34").Value = 44272
$ws.Range("J34").Value = 100
$ws.Range("D35").Value = 44272
$ws.Range("J35").Value = 50
$ws.Range("D36").Value = 44349
$ws.Range("J36").Value = 200
$ws.Range("K36").Value = 11000
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = 11500
$ws.Range("P36").Value = 639
$ws.Range("D37").Value = 44349
$ws.Range("J37").Value = 100
$ws.Range("K37").Value = 10000
$ws.Range("L37").Value = 10000
$ws.Range("M37").Value = 10000
$ws.Range("P37").Value = 556
$ws.Range("D38").Value = 44335
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 12000
$ws.Range("L38").Value = 13000
$ws.Range("M38").Value = 12500
$ws.Range("P38").Value = 694
$ws.Range("D39").Value = 44335
$ws.Range("J39").Value = 50
$ws.Range("K39").Value = 10000
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = 10000
$ws.Range("P39").Value = 556
